# Fixed citation graph with full text
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 (DOI 10.4103/ijo.IJO_1089_20) ---
# "Other found locations" cleared to a blank text cell (leading-apostrophe
# forces a text-typed empty string instead of a truly blank cell, then
# ClearFormats drops the quote-prefix style so no stray "s" attribute is left).
$ws.Cells.Item(2, 9).Value = "'"
$ws.Cells.Item(2, 9).ClearFormats()
# Authors: trailing row-version-id flag flipped 1 -> 0, padding whitespace normalised
$ws.Range("E2").Value = "[Anthony V%Das%NULL%0, Padmaja K%Rani%NULL%1, Pravin K%Vaddavalli%NULL%1]"
# "Misc. Data" now records the publisher
$ws.Range("J2").Value = "Wolters Kluwer - Medknow"

# --- Row 3 (DOI 10.1007/s40123-020-00269-3) ---
$ws.Cells.Item(3, 9).Value = "'"
$ws.Cells.Item(3, 9).ClearFormats()
# Abstract replaced with the full-text version (keeps the id="ParN"> markers)
$ws.Range("D3").Value = "Introduction`nid=`"Par1`">The coronavirus disease 2019 (COVID-19) pandemic has led to rapid adoption of teleophthalmology to deliver eyecare remotely.`n`n The purpose of our study was to assess the implementation and patient acceptability of video consultation for outpatient ophthalmic care at our institution.`n`n`nMethods`nid=`"Par2`">We conducted a retrospective, cross-sectional analysis and patient survey of adult patients who completed a virtual video visit at our institution from 18 March 18 through to 27 April 2020. All video visit encounters were assessed for patient characteristics, diagnoses, management, and follow-up outcomes.`n`n Patients were surveyed for their feedback on acceptability and utility of their virtual video consultation.`n`n`nResults`nid=`"Par3`">A total of 219 patients (mean age 55 years; range 21" + [char]0x2013 + "89 years) completed 231 video visit encounters at our department over a 6-week period, of whom 118 were women (54%).`n`n About half of these encounters were acute visits (102 visits, 47%).`n`n The most common diagnosis of these visits was postoperative state (20 visits, 9% of the total), followed by conjunctivitis (16 visits, 7%), and keratitis (14 visits, 6%).`n`n The most common management decisions were medication prescription (102 visits, 46%) or reassurance (86 visits, 39%), while 17 video visit patients (8%) were escalated to an urgent, in-person evaluation.`n`n Ninety-two patients completed a follow-up survey (42% response rate), of whom 45 (49%) indicated that they might have delayed seeking care during this pandemic in the absence of a virtual video option.`n`n Seventy-two (78%) reported that they would consider participating in a video visit as an alternative to an office-based encounter in the future, and the overall video visit experience was rated highly, with a weighted mean Likert scale rating of 4.3 out of 5 (Cronbach" + [char]0x2019 + "s " + [char]0x3B1 + [char]0x2009 + "=" + [char]0x2009 + "0.88).`n`n`nConclusion`nid=`"Par4`">Virtual video visits may be used to manage a range of ophthalmic complaints.`n`n Patients participating in this survey found such video visits acceptable and timesaving, and the majority would consider using video consultations for future eyecare encounters.`n`n`n"
# Re-run autofit so the much longer abstract does not leave an explicit custom
# row height behind (matches the source workbook, which has none on row 3).
$ws.Rows.Item(3).AutoFit()
# Authors: trailing row-version-id flag flipped 1 -> 0, padding whitespace normalised
$ws.Range("E3").Value = "[Gagan%Kalra%NULL%0, Andrew M.%Williams%NULL%1, Patrick W.%Commiskey%NULL%1, Eve M. R.%Bowers%NULL%1, Tadhg%Schempf%NULL%1, Jos" + [char]0xE9 + "-Alain%Sahel%NULL%1, Evan L.%Waxman%waxmane@upmc.edu%1, Roxana%Fu%fur3@upmc.edu%1]"
# "Misc. Data" now records the publisher
$ws.Range("J3").Value = "Springer Healthcare"

